{"js": "// Office.js (Word JavaScript API) script.\n// Applies two changes described by the diff:\n//  1. Adds a new bookmark \"OLE_LINK77\" around the word \"signsofleak\"\n//     (same span as the existing OLE_LINK72 bookmark there). Word\n//     renumbers every following bookmark's w:id by +1 automatically\n//     when the package is saved.\n//  2. Merges the two adjacent runs \"section.\" and \"eee\" (inside the\n//     \"Life Expectancy (EEE)\" row) into a single run \"section.eee\".\n\n// --- 1. Insert the OLE_LINK77 bookmark around \"signsofleak\" ---\nconst leakResults = context.document.body.search(\"signsofleak\", { matchCase: true, matchWholeWord: true });\nleakResults.load(\"text\");\nawait context.sync();\n\nif (leakResults.items.length > 0) {\n  leakResults.items[0].insertBookmark(\"OLE_LINK77\");\n  await context.sync();\n}\n\n// --- 2. Merge \"section.\" + \"eee\" into a single run \"section.eee\" ---\nconst eeeResults = context.document.body.search(\"section.eee\", { matchCase: true, matchWholeWord: false });\neeeResults.load(\"text\");\nawait context.sync();\n\nif (eeeResults.items.length > 0) {\n  eeeResults.items[0].insertText(\"section.eee\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies two changes described by the diff:\n#  1. Adds a new bookmark \"OLE_LINK77\" around the word \"signsofleak\"\n#     (same span as the existing OLE_LINK72 bookmark there). Word\n#     renumbers every following bookmark's w:id by +1 automatically\n#     when the package is saved.\n#  2. Merges the two adjacent runs \"section.\" and \"eee\" (inside the\n#     \"Life Expectancy (EEE)\" row) into a single run \"section.eee\".\n\n$d = $word.ActiveDocument\n\n# --- 1. Insert the OLE_LINK77 bookmark around \"signsofleak\" ---\n$find = $d.Content\n$find.Find.ClearFormatting()\n$find.Find.Text = \"signsofleak\"\n$find.Find.MatchCase = $true\n$find.Find.MatchWholeWord = $true\nif ($find.Find.Execute()) {\n    $d.Bookmarks.Add(\"OLE_LINK77\", $find)\n}\n\n# --- 2. Merge \"section.\" + \"eee\" into a single run \"section.eee\" ---\n$find2 = $d.Content\n$find2.Find.ClearFormatting()\n$find2.Find.Text = \"section.eee\"\n$find2.Find.MatchCase = $true\nif ($find2.Find.Execute()) {\n    # The target text is identical to the text already present (only the\n    # run split changes), so a direct re-assignment is a no-op. Round-trip\n    # through a temporary value to force Word to rewrite the range as a\n    # single run.\n    $find2.Text = \"__TMP_MERGE__\"\n    $find2.Text = \"section.eee\"\n}\n"}
